# Rename the 'Codelists' sheet to 'Cells' and make it the active sheet,
# with its selection moved to F19 (previously the 'Table' sheet was the
# active/selected one).

$wb = $excel.ActiveWorkbook

$wsCells = $wb.Worksheets.Item("Codelists")
$wsCells.Name = "Cells"

# Activating this sheet clears tabSelected on the previously active
# 'Table' sheet, sets it here instead, and updates the workbook's
# bookViews/workbookView activeTab accordingly.
$wsCells.Activate()

# Update the selection shown on this sheet.
$wsCells.Range("F19").Select() | Out-Null
